$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.597.85"
$ws.Range("E2").Value = "'  +4.67%  "
$ws.Range("D3").Value = "'3.612.79"
$ws.Range("E3").Value = "'  +4.62%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.17%  "
$ws.Range("D5").Value = "'630.00"
$ws.Range("E5").Value = "'  +4.84%  "
$ws.Range("D6").Value = "'158.74"
$ws.Range("E6").Value = "'  +8.46%  "
$ws.Range("D7").Value = "'3.610.94"
$ws.Range("E7").Value = "'  +4.70%  "
$ws.Range("E8").Value = "'  -0.16%  "
$ws.Range("D9").Value = "'0.495"
$ws.Range("E9").Value = "'  +4.43%  "
$ws.Range("D10").Value = "'0.150"
$ws.Range("E10").Value = "'  +10.80%  "
$ws.Range("D11").Value = "'7.49"
$ws.Range("E11").Value = "'  +9.53%  "
$ws.Range("D12").Value = "'0.443"
$ws.Range("E12").Value = "'  +6.23%  "
$ws.Range("E13").Value = "'  +6.55%  "
$ws.Range("D14").Value = "'33.72"
$ws.Range("E14").Value = "'  +9.08%  "
$ws.Range("D15").Value = "'4.223.54"
$ws.Range("E15").Value = "'  +4.53%  "
$ws.Range("B16").Value = "'WrappedEther"
$ws.Range("C16").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'3.611.72"
$ws.Range("E16").Value = "'  +4.60%  "
$ws.Range("B17").Value = "'WrappedBTC"
$ws.Range("C17").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "'69.491.55"
$ws.Range("E17").Value = "'  +4.42%  "
$ws.Range("E18").Value = "'  +1.15%  "
$ws.Range("D19").Value = "'6.75"
$ws.Range("E19").Value = "'  +7.39%  "
$ws.Range("D20").Value = "'16.22"
$ws.Range("E20").Value = "'  +9.80%  "
$ws.Range("D21").Value = "'10.28"
$ws.Range("E21").Value = "'  +15.23%  "
$ws.Range("D22").Value = "'462.43"
$ws.Range("E22").Value = "'  +5.78%  "
$ws.Range("D23").Value = "'0.647"
$ws.Range("E23").Value = "'  +4.79%  "
$ws.Range("D24").Value = "'79.01"
$ws.Range("E24").Value = "'  +3.15%  "
$ws.Range("D25").Value = "'0.0000135"
$ws.Range("E25").Value = "'  +10.67%  "
$ws.Range("D26").Value = "'10.75"
$ws.Range("E26").Value = "'  +7.97%  "
$ws.Range("D27").Value = "'3.757.02"
$ws.Range("E27").Value = "'  +4.52%  "
$ws.Range("E28").Value = "'  -0.01%  "
$ws.Range("D29").Value = "'9.40"
$ws.Range("E29").Value = "'  +15.20%  "
$ws.Range("E30").Value = "'  +6.70%  "
$ws.Range("D31").Value = "'1.73"
$ws.Range("E31").Value = "'  +13.76%  "
$ws.Range("D32").Value = "'0.173"
$ws.Range("E32").Value = "'  +9.04%  "
$ws.Range("D33").Value = "'6.57"
$ws.Range("E33").Value = "'  +8.51%  "
$ws.Range("E34").Value = "'  +0.11%  "
$ws.Range("B35").Value = "'ImmutableX"
$ws.Range("C35").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'1.96"
$ws.Range("E35").Value = "'  +7.22%  "
$ws.Range("B36").Value = "'EthereumClassic"
$ws.Range("C36").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "'26.57"
$ws.Range("E36").Value = "'  +5.00%  "
$ws.Range("D37").Value = "'3.609.93"
$ws.Range("E37").Value = "'  +4.91%  "
$ws.Range("E38").Value = "'  +9.09%  "
$ws.Range("D39").Value = "'2.43"
$ws.Range("E39").Value = "'  +15.11%  "
$ws.Range("E40").Value = "'  -0.03%  "
$ws.Range("D41").Value = "'0.0928"
$ws.Range("E41").Value = "'  +8.80%  "
$ws.Range("D42").Value = "'178.98"
$ws.Range("E42").Value = "'  +3.43%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "'  -0.01%  "
$ws.Range("D44").Value = "'5.69"
$ws.Range("E44").Value = "'  +6.20%  "
$ws.Range("D45").Value = "'31.88"
$ws.Range("E45").Value = "'  +24.17%  "
$ws.Range("E46").Value = "'  +4.69%  "
$ws.Range("E47").Value = "'  +15.21%  "
$ws.Range("B48").Value = "'dogwifhat"
$ws.Range("C48").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "'2.76"
$ws.Range("E48").Value = "'  +12.99%  "
$ws.Range("B49").Value = "'OKB"
$ws.Range("C49").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'45.95"
$ws.Range("E49").Value = "'  +1.68%  "
$ws.Range("D50").Value = "'7.85"
$ws.Range("E50").Value = "'  +5.01%  "
$ws.Range("D51").Value = "'0.269"
$ws.Range("E51").Value = "'  +10.97%  "
